$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This edit regenerates the localization-status report: the records for
# "e7874dfd-8dff-46b0-8b98-2c985fb33c4f" and "8e9ff772-6268-4785-9b07-524a7d25be06"
# swap places (row 4 <-> row 5 on every sheet), and the status of
# "6bc7d22b-d9a4-4551-afc2-dd7271add7aa" (row 3) and the relocated
# "e7874dfd..." record (now row 4) moves from "Ready for handoff" to
# "In Translation". The "8e9ff772..." record (now row 5) keeps
# "Ready for handoff".
# ---------------------------------------------------------------------------

# ---- Sheet "Overview" ----
$ws = $wb.Worksheets.Item("Overview")

# Row 3 (6bc7d22b...) status flips to "In Translation"
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

# Row 4 becomes the "e7874dfd..." record, status "In Translation"
$ws.Range("A4").Value = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("B4").Value = "e2e\e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"

# Row 5 becomes the "8e9ff772..." record, status stays "Ready for handoff"
$ws.Range("A5").Value = "8e9ff772-6268-4785-9b07-524a7d25be06.md"
$ws.Range("B5").Value = "e2e\8e9ff772-6268-4785-9b07-524a7d25be06.md"
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"

# Keep the hyperlink display text in sync with the new cell text (the
# r:id / target relationships themselves are unaffected).
$ws.Range("B4").Hyperlinks(1).TextToDisplay = "e2e\e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("B5").Hyperlinks(1).TextToDisplay = "e2e\8e9ff772-6268-4785-9b07-524a7d25be06.md"

# ---- Sheet "zh-cn" ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.377bfffacbbaae3ddd145dedada3fa890895f705.zh-cn.xlf"

$ws.Range("A5").Value = "8e9ff772-6268-4785-9b07-524a7d25be06.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("G5").Value = "8e9ff772-6268-4785-9b07-524a7d25be06.11ac84416452313dae1ffce2a6d6117bfef3413f.zh-cn.xlf"

$ws.Range("A4").Hyperlinks(1).TextToDisplay = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("A5").Hyperlinks(1).TextToDisplay = "8e9ff772-6268-4785-9b07-524a7d25be06.md"

# ---- Sheet "de-de" ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.377bfffacbbaae3ddd145dedada3fa890895f705.de-de.xlf"

$ws.Range("A5").Value = "8e9ff772-6268-4785-9b07-524a7d25be06.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("G5").Value = "8e9ff772-6268-4785-9b07-524a7d25be06.11ac84416452313dae1ffce2a6d6117bfef3413f.de-de.xlf"

$ws.Range("A4").Hyperlinks(1).TextToDisplay = "e7874dfd-8dff-46b0-8b98-2c985fb33c4f.md"
$ws.Range("A5").Hyperlinks(1).TextToDisplay = "8e9ff772-6268-4785-9b07-524a7d25be06.md"
